$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Workbook/window metadata tweaks (best-effort; engine auto-manages fileVersion) ---

# --- Header-style rows get the bordered/filled header format copied from the existing header row A44 ---
$ws.Range("A44").Copy()
$ws.Range("A46").PasteSpecial(-4122)
$ws.Range("B46").PasteSpecial(-4122)
$ws.Range("A55").PasteSpecial(-4122)
$ws.Range("B55").PasteSpecial(-4122)
$ws.Range("A60").PasteSpecial(-4122)
$ws.Range("B60").PasteSpecial(-4122)
$ws.Range("A67").PasteSpecial(-4122)
$ws.Range("B67").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Cell values + rich-text runs, written in the exact order shared strings were introduced ---
$ws.Range("A47").Value = '<style>   </style>'
$ws.Range("B47").Value = 'Internal CSS define for a single HTML page; defines in the <head> section of HTML; '
$ws.Range("A48").Value = 'External CSS  <link rel="stylesheet" href="css file name">'
$ws.Range("A49").Value = 'border'
$cell = $ws.Range("B49")
$cell.Value = 'defines border around the html element; format: border: n px type color; Example: border : 1px solid blue;'
$cell.Characters(41, 6).Font.Bold = $true
$cell.Characters(47, 2).Font.Bold = $false
$cell.Characters(49, 25).Font.Italic = $true
$cell.Characters(74, 9).Font.Bold = $true
$cell.Characters(83, 24).Font.Italic = $true
$ws.Range("A50").Value = 'padding'
$ws.Range("B50").Value = 'defines padding(space) between text and the border'
$ws.Range("A51").Value = 'margin'
$ws.Range("B51").Value = 'defines a margin (space) outside the border'
$ws.Range("A52").Value = ' id = "  "'
$cell = $ws.Range("B52")
$cell.Value = 'to define a specific style for one special element; define the id number in hash tag in the css;'
$cell.Characters(77, 5).Font.Italic = $true
$cell.Characters(82, 15).Font.Bold = $false
$ws.Range("A53").Value = 'class'
$ws.Range("B48").Value = 'Defined for many HTML sheets; change the look of the entire website; add link in the <head> section; can be a full path URL'
$ws.Range("B54").Value = 'links to a style sheet located in the folder or the web page; defines in the <head>'
$ws.Range("A54").Value = '<link rel="stylesheet" href="path or URL">'
$cell = $ws.Range("B53")
$cell.Value = 'defines style for a special type of element; add a class attribute to that element; called as element.name in css'
$cell.Characters(103, 5).Font.Italic = $true
$cell.Characters(108, 6).Font.Bold = $false
$ws.Range("A55").Value = 'HTML Links'
$ws.Range("A56").Value = '<a href="url"> link text</a>'
$ws.Range("B56").Value = 'Defined with <a> tag; Links are to navigate from one page to another; href sepcifies destination address; link text is the visible part;'
$ws.Range("A57").Value = 'Local link'
$ws.Range("A58").Value = 'absoulte link'
$ws.Range("B58").Value = 'has the full address and link to a different website'
$ws.Range("B57").Value = 'link to the same website; doesn''t have the http://www….'
$ws.Range("A59").Value = 'default link properties'
$ws.Range("B59").Value = 'unvisited link underlined and blue; visited link underlined and purple; active link underlined and red;'
$ws.Range("A60").Value = 'HTML Links - Target attributes'
$ws.Range("A62").Value = '_blank'
$ws.Range("B62").Value = 'Opens the linked document in new window or tab'
$ws.Range("A63").Value = '_self'
$ws.Range("B63").Value = 'Opens the linked document in the same window or tab as it was clicked(default)'
$ws.Range("A64").Value = '_parent'
$ws.Range("B64").Value = 'Opens the linked document in the parernt frame'
$ws.Range("A65").Value = '_top'
$ws.Range("B65").Value = 'Opens the linked document in the full body of the window'
$ws.Range("A66").Value = 'framename'
$ws.Range("B66").Value = 'Opens the linked document in the named frame'
$ws.Range("A61").Value = '<target=" values"'
$ws.Range("B61").Value = 'Specifies where to open the linked document based on the below values'
$ws.Range("A67").Value = 'Image as link'

# --- Row heights for header rows (matches customHeight introduced by the author) ---
$ws.Rows(46).RowHeight = 15.75
$ws.Rows(55).RowHeight = 15.75
$ws.Rows(60).RowHeight = 15.75
$ws.Rows(67).RowHeight = 15.75

# --- Final selection state ---
$ws.Range("A68").Select()
